$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates - B2 and D2 get new values, C2 and E2 become empty
$ws.Range("B2").Value = 23.254962237594334
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 34.561854057171395
$ws.Range("E2").ClearContents()

# Row 3 updates
$ws.Range("B3").Value = 20.682618561610106
$ws.Range("C3").Value = -6.5016201590062561
$ws.Range("D3").Value = 23.677839492541434
$ws.Range("E3").Value = -12.316003057273068

# Update the selection to reflect the new active range
$ws.Range("B1:E3").Select()
